$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Version 0.1.6 -> 0.1.7
$wsMeta.Range("B3").Value = "0.1.7"

# Status active -> draft
$wsMeta.Range("B6").Value = "draft"

# Date updated
$wsMeta.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# Contact (publisher contact) updated
$wsMeta.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11 already held a second "Contact" row (duplicate of the old display text) -
# repurpose it for the individual contact.
$wsMeta.Range("A11").Value = "Contact"
$wsMeta.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row at 12 for Jurisdiction, pushing Description/Purpose/Copyright/Immutable down.
$wsMeta.Rows("12:12").Insert()

$wsMeta.Range("A12").Value = "Jurisdiction"
$wsMeta.Range("B12").Value = ""

# Apply the same formatting (border/fill/alignment) as the other data rows to the new row.
$wsMeta.Range("A11:B11").Copy()
$wsMeta.Range("A12:B12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
